$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Vaccine_Uptake (column D) values per state, keyed by row number.
$updates = @{
    2 = 0.67789221829395441
    3 = 0.63662190843520516
    4 = 0.65763096820159372
    5 = 0.66981387512608215
    6 = 0.66981309319895665
    7 = 0.66994337523819347
    8 = 0.72233867299461274
    9 = 0.6700321169479786
    10 = 0.6669743959519292
    11 = 0.69619952087491499
    12 = 0.66737853952129467
    13 = 0.69412808752690369
    15 = 0.70181120458206547
    16 = 0.6785643000554713
    17 = 0.66355568254990338
    18 = 0.70047026737419282
    19 = 0.68524113843146095
    20 = 0.68291066818134882
    21 = 0.71296602697095435
    22 = 0.73121809200352961
    23 = 0.72692201185717209
    24 = 0.71849866060023149
    25 = 0.69043666416593508
    26 = 0.67749262182566916
    27 = 0.69667395526471732
    28 = 0.64393686709965769
    29 = 0.7098582699261563
    30 = 0.66625236937255883
    31 = 0.73555820308754438
    32 = 0.65147372382932933
    33 = 0.67892409930004327
    34 = 0.67328727911957031
    35 = 0.72185349483294137
    36 = 0.6720543500222248
    37 = 0.69746488696925146
    38 = 0.68862360939683165
    39 = 0.6951508413994113
    40 = 0.67936949300451122
    42 = 0.67749962620459414
    43 = 0.69890255760248643
    44 = 0.69763079319839627
    45 = 0.71067323586547015
    46 = 0.6990675430858273
    47 = 0.73104171685098007
    48 = 0.68884856819917761
    49 = 0.70203227943135027
    50 = 0.65276166386902434
    51 = 0.70714155816654589
    52 = 0.69098971631205675
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 4).Value = $updates[$row]
}
